# Workbook/worksheet handles (per task: $wb already resolves to the open
# workbook, but we re-derive it from $excel to be safe/explicit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the column headers in row 1 so the "_old"/"_new" suffixes
#    become "_FV2304"/"_FV2310" (the two EDIFACT format-version labels
#    this diff/merge compares), leaving "diff" untouched.
# ---------------------------------------------------------------------
$lastCol = 20
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Text
    if ($current -like "*_old") {
        $cell.Value = ($current -replace "_old$", "_FV2304")
    } elseif ($current -like "*_new") {
        $cell.Value = ($current -replace "_new$", "_FV2310")
    }
}

# ---------------------------------------------------------------------
# 2) Turn the used range A1:T72 into an Excel Table ("Table1") that
#    picks up the renamed headers, with a normal auto-filter and
#    default banded-row styling.
# ---------------------------------------------------------------------
$tableRange = $ws.Range("A1:T72")
$tbl = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------
# 3) Freeze the header row: split/freeze after row 1, so row 2 becomes
#    the top-left cell of the scrolling pane.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
